# "attdence 23-Jan-24 and 24-Jan-24"
# Add two new attendance-taking dates (23-Jan-2024 and 24-Jan-2024) as
# columns Y and Z on the "Jan-2024" sheet, mirroring each student's
# attendance mark from the most recent existing date (column W, 21-Jan-24).
# Also touch the four still-unused trailing columns AC:AF so the sheet's
# used range/formatting/data-validation grow the same way they did in the
# authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jan-2024")

# --- Header row: stamp the two new session dates ---------------------
$ws.Range("Y1").Value = 45314          # 23-Jan-2024
$ws.Range("Y1").NumberFormat = "d-mmm-yy"
$ws.Range("Z1").Value = 45315          # 24-Jan-2024
$ws.Range("Z1").NumberFormat = "d-mmm-yy"

# --- Attendance rows: same mark as the last recorded date (col W) -----
$ws.Range("Y2").Value = "Present"
$ws.Range("Z2").Value = "Present"

$ws.Range("Y3").Value = "Present"
$ws.Range("Z3").Value = "Present"

$ws.Range("Y4").Value = "Absent"
$ws.Range("Z4").Value = "Absent"

$ws.Range("Y5").Value = "Absent"
$ws.Range("Z5").Value = "Absent"

$ws.Range("Y6").Value = "Absent"
$ws.Range("Z6").Value = "Absent"

$ws.Range("Y7").Value = "Present"
$ws.Range("Z7").Value = "Present"

$ws.Range("Y8").Value = "Absent"
$ws.Range("Z8").Value = "Absent"

# --- Bring the new trailing blank columns AC:AF into the used grid, ---
# --- copying the formatting of the adjacent column AB so no new -------
# --- cell styles are introduced. ---------------------------------------
$ws.Range("AB1").Copy()
$ws.Range("AC1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AB2:AB8").Copy()
$ws.Range("AC2:AF8").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Extend the "Present, Absent, Reason" list validation over the ----
# --- newly-used columns -------------------------------------------------
$ws.Range("C2:AF8").Validation.Delete()
$ws.Range("C2:AF8").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')

# --- Update the on-screen selection to the last cell touched ----------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 21
$ws.Range("Z8").Select()
